# Auto-generated script applying the coinranking data refresh described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.971.93"
$ws.Range("E2").Value = "  -5.79%  "
$ws.Range("D3").Value = "3.250.41"
$ws.Range("E3").Value = "  -7.00%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "177.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -11.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "514.36"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.597"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.71%  "
$ws.Range("D8").Value = "3.244.84"
$ws.Range("E8").Value = "  -6.87%  "
$ws.Range("E9").Value = "  -0.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.616"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.63%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "58.05"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.66%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.132"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -7.80%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000255"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -6.81%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.07"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.92%  "
$ws.Range("D15").Value = "3.761.22"
$ws.Range("E15").Value = "  -7.29%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.117"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.97%  "
$ws.Range("D17").Value = "3.240.50"
$ws.Range("E17").Value = "  -7.49%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.44"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.79%  "
$ws.Range("D19").Value = "62.808.11"
$ws.Range("E19").Value = "  -5.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.95"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.83%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.948"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -7.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "370.41"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.34%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.14"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.44%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.46%  "
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.68"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -7.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.86"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.06"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.41"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.75%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.63"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.30"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.94%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "28.48"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.75%  "
$ws.Range("E32").Value = "  -6.40%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "631.12"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.39%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.27"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.25%  "
$ws.Range("E35").Value = "  -3.83%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "58.56"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.39%  "
$ws.Range("E37").Value = "  -1.08%  "
$ws.Range("E38").Value = "  +0.14%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.41"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.86%  "
$ws.Range("E40").Value = "  -0.06%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.124"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.17%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "2.922.39"
$ws.Range("E42").Value = "  -5.89%  "
$ws.Range("D43").Value = "0.0$([char]8323)0662"
$ws.Range("E43").Value = "  -5.45%  "
$ws.Range("B44").Value = "ThetaToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.68"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -11.86%  "
$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.42"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.49%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0392"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.80%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.84"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +9.34%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.60"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.51%  "
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.125"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.13%  "
$ws.Range("B50").Value = "ApeXProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.95"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.94%  "
$ws.Range("E51").Value = "  -11.00%  "
